# Weekly crime-data refresh for the 123rd Precinct CompStat sheet.
# Updates: report header (issue number, week-of dates) and the
# Week-to-Date / 28-Day / Year-to-Date / 2-Year complaint tables
# (rows 14-29) with newly collected figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Report header: issue number and covered-week dates ---
# "Volume 31   Number  6"  ->  "Volume 31   Number  7"
$ws.Range("A8").Characters(21, 1).Text = "7"

# "Report Covering the Week  2/5/2024  Through  2/11/2024"
# -> "Report Covering the Week  2/12/2024  Through  2/18/2024"
$ws.Range("C9").Characters(27, 8).Text = "2/12/2024"
$ws.Range("C9").Characters(47, 9).Text = "2/18/2024"


# --- Simple value-only updates (style/type unchanged) ---
$ws.Range("F16").Value = 3
$ws.Range("I16").Value = 3
$ws.Range("J16").Value = 3
$ws.Range("K16").Value = 0
$ws.Range("M16").Value = 200
$ws.Range("N16").Value = -50
$ws.Range("C17").Value = 4
$ws.Range("D17").Value = 2
$ws.Range("E17").Value = 100
$ws.Range("F17").Value = 7
$ws.Range("G17").Value = 5
$ws.Range("H17").Value = 40
$ws.Range("I17").Value = 10
$ws.Range("J17").Value = 9
$ws.Range("K17").Value = 11.111111111111
$ws.Range("L17").Value = 11.111111111111
$ws.Range("M17").Value = 42.857142857142
$ws.Range("N17").Value = 66.666666666666
$ws.Range("L18").Value = 50
$ws.Range("M18").Value = -78.571428571428
$ws.Range("N18").Value = -92.5
$ws.Range("C19").Value = 5
$ws.Range("E19").Value = 0
$ws.Range("F19").Value = 23
$ws.Range("G19").Value = 21
$ws.Range("H19").Value = 9.523809523809
$ws.Range("I19").Value = 29
$ws.Range("J19").Value = 34
$ws.Range("K19").Value = -14.705882352941
$ws.Range("L19").Value = -25.641025641025
$ws.Range("M19").Value = 20.833333333333
$ws.Range("N19").Value = 81.25
$ws.Range("F20").Value = 2
$ws.Range("G20").Value = 2
$ws.Range("L20").Value = -64.705882352941
$ws.Range("M20").Value = -33.333333333333
$ws.Range("N20").Value = -93.617021276595
$ws.Range("C21").Value = 10
$ws.Range("D21").Value = 9
$ws.Range("E21").Value = 11.111111111111
$ws.Range("F21").Value = 37
$ws.Range("G21").Value = 34
$ws.Range("H21").Value = 8.823529411764
$ws.Range("I21").Value = 54
$ws.Range("J21").Value = 65
$ws.Range("K21").Value = -16.923076923076
$ws.Range("L21").Value = -19.402985074626
$ws.Range("M21").Value = -1.818181818181
$ws.Range("N21").Value = -66.871165644171
$ws.Range("C24").Value = 7
$ws.Range("D24").Value = 11
$ws.Range("E24").Value = -36.363636363636
$ws.Range("F24").Value = 29
$ws.Range("H24").Value = -30.952380952381
$ws.Range("I24").Value = 53
$ws.Range("J24").Value = 77
$ws.Range("K24").Value = -31.168831168831
$ws.Range("L24").Value = 10.416666666666
$ws.Range("M24").Value = -27.397260273972
$ws.Range("C25").Value = 7
$ws.Range("D25").Value = 8
$ws.Range("E25").Value = -12.5
$ws.Range("F25").Value = 13
$ws.Range("G25").Value = 17
$ws.Range("H25").Value = -23.529411764705
$ws.Range("I25").Value = 18
$ws.Range("J25").Value = 27
$ws.Range("K25").Value = -33.333333333333
$ws.Range("L25").Value = -21.739130434782
$ws.Range("M25").Value = -41.935483870967
$ws.Range("F26").Value = 2
$ws.Range("G27").Value = 1

# --- Updates requiring a style/type change (text<->number) ---
$ws.Range("F14").NumberFormat = "@"
$ws.Range("F14").Value = "0"
$ws.Range("A14").Copy()
$ws.Range("F14").PasteSpecial(-4122)

$ws.Range("N14").Value = 0
$ws.Range("H18").Copy()
$ws.Range("N14").PasteSpecial(-4122)

$ws.Range("F15").NumberFormat = "@"
$ws.Range("F15").Value = "0"
$ws.Range("A14").Copy()
$ws.Range("F15").PasteSpecial(-4122)

$ws.Range("D16").Value = 2
$ws.Range("I14").Copy()
$ws.Range("D16").PasteSpecial(-4122)

$ws.Range("E16").Value = -50
$ws.Range("H18").Copy()
$ws.Range("E16").PasteSpecial(-4122)

$ws.Range("G16").Value = 2
$ws.Range("I14").Copy()
$ws.Range("G16").PasteSpecial(-4122)

$ws.Range("H16").Value = 50
$ws.Range("H18").Copy()
$ws.Range("H16").PasteSpecial(-4122)

$ws.Range("C18").NumberFormat = "@"
$ws.Range("C18").Value = "0"
$ws.Range("A14").Copy()
$ws.Range("C18").PasteSpecial(-4122)

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0"
$ws.Range("A14").Copy()
$ws.Range("D18").PasteSpecial(-4122)

$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "***.*"
$ws.Range("A14").Copy()
$ws.Range("E18").PasteSpecial(-4122)

$ws.Range("C20").NumberFormat = "@"
$ws.Range("C20").Value = "0"
$ws.Range("A14").Copy()
$ws.Range("C20").PasteSpecial(-4122)

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0"
$ws.Range("A14").Copy()
$ws.Range("D27").PasteSpecial(-4122)

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "***.*"
$ws.Range("A14").Copy()
$ws.Range("E27").PasteSpecial(-4122)

$ws.Range("N28").Value = -100
$ws.Range("H18").Copy()
$ws.Range("N28").PasteSpecial(-4122)

$ws.Range("N29").Value = -100
$ws.Range("H18").Copy()
$ws.Range("N29").PasteSpecial(-4122)
